$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New job numbers (B2:B4) must stay text, matching how the existing
# sheet stores them as shared strings rather than numbers.
$ws.Range("B2").Value = "'32255213"
$ws.Range("B3").Value = "'32255214"
$ws.Range("B4").Value = "'32255215"

# Restore the default (unstyled) cell format so Excel's "stored as text"
# warning formatting doesn't introduce a new style index.
$ws.Range("B2:B4").Style = "Normal"
